$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$sched = $wb.Worksheets.Item("Schedule")
$sched.Range("E2").Value = 356.2560397499999
$sched.Range("F2").Value = 7.853969130291003
$sched.Range("E3").Value = 430.923909
$sched.Range("F3").Value = 28.50025853174603

# --- Sheet "Detailed" ---
$det = $wb.Worksheets.Item("Detailed")
$det.Range("B12").Value = 67.38898
$det.Range("B13").Value = 81.43025
$det.Range("B14").Value = 78
$det.Range("C14").Value = "historical"
$det.Range("B15").Value = 59.33406
$det.Range("B19").Value = -6.21646
$det.Range("B20").Value = -7.81926
$det.Range("B21").Value = -7.94282
$det.Range("B22").Value = -7.49537
$det.Range("B23").Value = -7.48969
$det.Range("B24").Value = -7.54002
$det.Range("B25").Value = -5.50985
$det.Range("B26").Value = -6.74852
$det.Range("B27").Value = -7.97656
$det.Range("B28").Value = -7.8587
$det.Range("B29").Value = -6.07252
$det.Range("B34").Value = -9.782209999999999
$det.Range("B35").Value = -10
$det.Range("B37").Value = -7.34615
$det.Range("B38").Value = -0.39869
$det.Range("B39").Value = 7.27483
$det.Range("B40").Value = 29.36849
$det.Range("B41").Value = 55.33036
$det.Range("B43").Value = 53.24127
$det.Range("B46").Value = 54.12603
$det.Range("B49").Value = 52.66218
